# Edit RAF for hydrogen power plants and enable CCS power plants for reliability pass
$wb = $excel.ActiveWorkbook

# --- ESUfRaLCD-reliability sheet: insert 4 rows for CCS power-plant types -------------
$wsRel = $wb.Worksheets.Item("ESUfRaLCD-reliability")

# Insert 4 blank rows above row 8 (pushes hydrogen rows + blanks down by 4,
# matching Excel's own formula auto-fill behaviour on row insert).
$wsRel.Range("A8:A11").EntireRow.Insert()

# Column B first for the two brand-new "w ccs es" sources ...
$wsRel.Range("B8").Value = "hard coal w ccs es"
$wsRel.Range("B9").Value = "natural gas combined cycle w ccs es"

# ... then column A for those same two rows (new "w ccs" strings, no "es" suffix) ...
$wsRel.Range("A8").Value = "hard coal w ccs"
$wsRel.Range("A9").Value = "natural gas combined cycle w ccs"

# ... then column A for the biomass/lignite rows (new "w ccs" strings, no "es" suffix) ...
$wsRel.Range("A10").Value = "biomass w ccs"
$wsRel.Range("A11").Value = "lignite w ccs"

# ... and column B for those rows reuses the existing "w CCS es" strings.
$wsRel.Range("B10").Value = "biomass w CCS es"
$wsRel.Range("B11").Value = "lignite w CCS es"

# Column C: power-plant-name formula, same pattern as the rest of the column.
$wsRel.Range("C8").Formula = "=IF(A8=`"`",`"`",CONCATENATE(A8,`" power plants`"))"
$wsRel.Range("C9").Formula = "=IF(A9=`"`",`"`",CONCATENATE(A9,`" power plants`"))"
$wsRel.Range("C10").Formula = "=IF(A10=`"`",`"`",CONCATENATE(A10,`" power plants`"))"
$wsRel.Range("C11").Formula = "=IF(A11=`"`",`"`",CONCATENATE(A11,`" power plants`"))"

# --- About sheet: its selection moved too, but focus ultimately rests on reliability --
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("A29:A30").Select()

# Re-select a cell on the reliability sheet (matches post-edit selection snapshot) and
# leave it as the active/visible tab.
$wsRel.Activate()
$wsRel.Range("A12").Select()
